# Apply the data updates described in the commit "Write data org slides".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "ts" measurements for rows 3-5.
$ws.Range("F3").Value = 102.45
$ws.Range("F4").Value = 57.02
$ws.Range("F5").Value = 27.03

# Rows 17-21: reclassify system/location from
# "pit latrine"/"household" to "septic tank"/"public toilet".
$ws.Range("C17:C21").Value = "septic tank"
$ws.Range("D17:D21").Value = "public toilet"
